# Add data for 2025-12-26
# Updates year-to-date crime counts (as of Dec 26) across the
# Citywide Totals sheet, the By Neighborhood summary sheet, and
# the per-neighborhood detail sheets, reflecting newly published
# and backfilled Chicago Police Department crime records.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H3").Value = 161
$ws.Range("I3").Value = 198
$ws.Range("K3").Value = 229
$ws.Range("B9").Value = 390
$ws.Range("C9").Value = 503
$ws.Range("D9").Value = 435
$ws.Range("F9").Value = 574
$ws.Range("H9").Value = 474
$ws.Range("K9").Value = 537
$ws.Range("L9").Value = 455
$ws.Range("B10").Value = 1406
$ws.Range("C10").Value = 1659
$ws.Range("D10").Value = 1881
$ws.Range("E10").Value = 2303
$ws.Range("F10").Value = 2193
$ws.Range("H10").Value = 634
$ws.Range("I10").Value = 875
$ws.Range("K10").Value = 708
$ws.Range("L10").Value = 692
$ws.Range("B11").Value = 1939
$ws.Range("C11").Value = 2326
$ws.Range("D11").Value = 2563
$ws.Range("E11").Value = 3049
$ws.Range("F11").Value = 3022
$ws.Range("H11").Value = 1403
$ws.Range("I11").Value = 1745
$ws.Range("K11").Value = 1659

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 23

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D8").Value = 48
$ws.Range("D9").Value = 97

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K9").Value = 20
$ws.Range("K10").Value = 56

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L8").Value = 44
$ws.Range("L10").Value = 96

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 30
$ws.Range("B9").Value = 230
$ws.Range("E9").Value = 691
$ws.Range("F9").Value = 560
$ws.Range("H9").Value = 114
$ws.Range("B10").Value = 279
$ws.Range("E10").Value = 781
$ws.Range("F10").Value = 647
$ws.Range("H10").Value = 233
$ws.Range("I10").Value = 322

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F5").Value = 24
$ws.Range("F7").Value = 59

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("E7").Value = 43
$ws.Range("E8").Value = 62

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B8").Value = 78
$ws.Range("H8").Value = 102
$ws.Range("K19").Value = 56
$ws.Range("C20").Value = 23
$ws.Range("D27").Value = 29
$ws.Range("D32").Value = 97
$ws.Range("B35").Value = 23
$ws.Range("K35").Value = 15
$ws.Range("L36").Value = 96
$ws.Range("F50").Value = 59
$ws.Range("B53").Value = 279
$ws.Range("E53").Value = 781
$ws.Range("F53").Value = 647
$ws.Range("H53").Value = 233
$ws.Range("I53").Value = 322
$ws.Range("L55").Value = 2
$ws.Range("E65").Value = 62
$ws.Range("H67").Value = 11
$ws.Range("H68").Value = 8
$ws.Range("D72").Value = 14
$ws.Range("F76").Value = 68
$ws.Range("K77").Value = 66
$ws.Range("C78").Value = 37
$ws.Range("K84").Value = 4
$ws.Range("F95").Value = 67
$ws.Range("I96").Value = 17
$ws.Range("B99").Value = 1939
$ws.Range("C99").Value = 2326
$ws.Range("D99").Value = 2563
$ws.Range("E99").Value = 3049
$ws.Range("F99").Value = 3022
$ws.Range("H99").Value = 1403
$ws.Range("I99").Value = 1745
$ws.Range("K99").Value = 1659

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("C5").Value = 33
$ws.Range("C6").Value = 37

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("F9").Value = 48
$ws.Range("F10").Value = 68

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("H6").Value = 6
$ws.Range("H8").Value = 11

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("B5").Value = 7
$ws.Range("K5").Value = 6
$ws.Range("B7").Value = 23
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("D7").Value = 20
$ws.Range("D8").Value = 29

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 7
$ws.Range("K8").Value = 19
$ws.Range("K10").Value = 66

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("F6").Value = 59
$ws.Range("F7").Value = 67

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I8").Value = 11
$ws.Range("I9").Value = 17

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("D4").Value = 6
$ws.Range("D6").Value = 14

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 2

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("H3").Value = 2
$ws.Range("H8").Value = 8

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("B8").Value = 23
$ws.Range("H9").Value = 41
$ws.Range("B10").Value = 78
$ws.Range("H10").Value = 102
